# Atualização automática de ALEGRETE.xlsx
#
# 1) Rename "Paineis DARQ"            -> "PAINEIS DARQ"
# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3) Delete the "Desarquivamentos Pendentes" sheet entirely
#    (the "DGC" sheet and all other sheets keep their data untouched)

$wb = $excel.ActiveWorkbook

$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
